{"js": "const replacements = [\n  [\"170\u00f79=\", \"770\u00f74=\"],\n  [\"981\u00f79=\", \"591\u00f76=\"],\n  [\"753\u00f72=\", \"524\u00f74=\"],\n  [\"948\u00f75=\", \"864\u00f72=\"],\n  [\"820\u00f72=\", \"834\u00f77=\"],\n  [\"374\u00f75=\", \"305\u00f72=\"],\n  [\"798\u00f77=\", \"371\u00f79=\"],\n  [\"627\u00f74=\", \"975\u00f74=\"],\n  [\"979\u00f78=\", \"735\u00f72=\"],\n  [\"322\u00f78=\", \"864\u00f74=\"],\n  [\"755\u00f73=\", \"674\u00f77=\"],\n  [\"328\u00f76=\", \"480\u00f78=\"],\n  [\"253\u00f77=\", \"830\u00f73=\"],\n  [\"835\u00f73=\", \"723\u00f72=\"],\n  [\"854\u00f77=\", \"165\u00f79=\"],\n  [\"604\u00f74=\", \"225\u00f73=\"],\n  [\"917\u00f78=\", \"693\u00f79=\"],\n  [\"710\u00f77=\", \"144\u00f78=\"],\n  [\"382\u00f75=\", \"209\u00f77=\"],\n  [\"790\u00f76=\", \"232\u00f74=\"],\n  [\"252\u00f76=\", \"599\u00f76=\"],\n  [\"431\u00f77=\", \"759\u00f72=\"],\n  [\"908\u00f78=\", \"842\u00f79=\"],\n  [\"492\u00f73=\", \"911\u00f77=\"],\n  [\"422\u00f77=\", \"206\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"170\u00f79=\", \"770\u00f74=\"),\n    @(\"981\u00f79=\", \"591\u00f76=\"),\n    @(\"753\u00f72=\", \"524\u00f74=\"),\n    @(\"948\u00f75=\", \"864\u00f72=\"),\n    @(\"820\u00f72=\", \"834\u00f77=\"),\n    @(\"374\u00f75=\", \"305\u00f72=\"),\n    @(\"798\u00f77=\", \"371\u00f79=\"),\n    @(\"627\u00f74=\", \"975\u00f74=\"),\n    @(\"979\u00f78=\", \"735\u00f72=\"),\n    @(\"322\u00f78=\", \"864\u00f74=\"),\n    @(\"755\u00f73=\", \"674\u00f77=\"),\n    @(\"328\u00f76=\", \"480\u00f78=\"),\n    @(\"253\u00f77=\", \"830\u00f73=\"),\n    @(\"835\u00f73=\", \"723\u00f72=\"),\n    @(\"854\u00f77=\", \"165\u00f79=\"),\n    @(\"604\u00f74=\", \"225\u00f73=\"),\n    @(\"917\u00f78=\", \"693\u00f79=\"),\n    @(\"710\u00f77=\", \"144\u00f78=\"),\n    @(\"382\u00f75=\", \"209\u00f77=\"),\n    @(\"790\u00f76=\", \"232\u00f74=\"),\n    @(\"252\u00f76=\", \"599\u00f76=\"),\n    @(\"431\u00f77=\", \"759\u00f72=\"),\n    @(\"908\u00f78=\", \"842\u00f79=\"),\n    @(\"492\u00f73=\", \"911\u00f77=\"),\n    @(\"422\u00f77=\", \"206\u00f76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}"}
